$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 230 (shifts existing rows 230..319 down to 231..320)
$ws.Rows.Item(230).Insert()

# Populate the newly inserted row 230 with the new record's data
$ws.Cells.Item(230, 1).Value = 10
$ws.Cells.Item(230, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(230, 3).Value = "La Araucanía"
$ws.Cells.Item(230, 4).Value = 45146
$ws.Cells.Item(230, 5).Value = 9
$ws.Cells.Item(230, 6).Value = "Fruta"
$ws.Cells.Item(230, 7).Value = 100104
$ws.Cells.Item(230, 8).Value = "Frutos de pepita"
$ws.Cells.Item(230, 9).Value = 100104003
$ws.Cells.Item(230, 10).Value = "Membrillo"
$ws.Cells.Item(230, 11).Value = "Champion"
$ws.Cells.Item(230, 12).Value = "Primera"
$ws.Cells.Item(230, 13).Value = 170
$ws.Cells.Item(230, 14).Value = 15000
$ws.Cells.Item(230, 15).Value = 16000
$ws.Cells.Item(230, 16).Value = 15471
$ws.Cells.Item(230, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(230, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(230, 19).Value = 860
$ws.Cells.Item(230, 20).Value = 18
